# Auto-generated edit script applying the cryptos.xlsx diff
# (GitHub Actions "Updated cryptos list" data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value is numeric-looking text (e.g. "1.17", "239.40").
# Excel auto-converts such strings to numbers on assignment (General format),
# which would also silently drop significant trailing zeros (e.g. 239.40 -> 239.4).
# Force Text format on just these cells first so the literal string is preserved.
$textForceCells = @(
    "D5", "D6", "D8", "D9", "D10", "D11", "D13", "D15", "D17", "D19", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D30", "D32", "D33", "D34", "D35", "D36", "D39", "D41", "D42", "D43", "D44", "D45", "D47", "D48", "D50"
)
foreach ($cellName in $textForceCells) {
    $ws.Range($cellName).NumberFormat = "@"
}

# Apply all updated cell values (row by row, matching the source diff).
$ws.Range("D2").Value = "36.530.62"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "2.098.86"
$ws.Range("E3").Value = "  +9.50%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "252.29"
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("D6").Value = "0.656"
$ws.Range("E6").Value = "  -6.44%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "47.42"
$ws.Range("E8").Value = "  +5.72%  "
$ws.Range("D9").Value = "59.94"
$ws.Range("E9").Value = "  +2.51%  "
$ws.Range("D10").Value = "0.375"
$ws.Range("E10").Value = "  +1.08%  "
$ws.Range("D11").Value = "0.0742"
$ws.Range("E11").Value = "  -2.94%  "
$ws.Range("E12").Value = "  -0.28%  "
$ws.Range("D13").Value = "14.59"
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("D14").Value = "2.401.15"
$ws.Range("E14").Value = "  +9.44%  "
$ws.Range("D15").Value = "0.824"
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("D16").Value = "2.093.95"
$ws.Range("E16").Value = "  +9.29%  "
$ws.Range("D17").Value = "5.07"
$ws.Range("E17").Value = "  -1.24%  "
$ws.Range("D18").Value = "36.503.54"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").Value = "72.68"
$ws.Range("E19").Value = "  -2.60%  "
$ws.Range("D20").Value = "0.0₃0826"
$ws.Range("E20").Value = "  -4.38%  "
$ws.Range("D21").Value = "13.16"
$ws.Range("E21").Value = "  -1.89%  "
$ws.Range("D22").Value = "239.40"
$ws.Range("E22").Value = "  -4.47%  "
$ws.Range("D23").Value = "5.15"
$ws.Range("E23").Value = "  -1.45%  "
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").Value = "2.45"
$ws.Range("E25").Value = "  -6.78%  "
$ws.Range("D26").Value = "170.20"
$ws.Range("E26").Value = "  +0.82%  "
$ws.Range("D27").Value = "21.40"
$ws.Range("E27").Value = "  +14.14%  "
$ws.Range("D28").Value = "9.10"
$ws.Range("E28").Value = "  +3.43%  "
$ws.Range("E29").Value = "  -10.24%  "
$ws.Range("D30").Value = "28.13"
$ws.Range("E30").Value = "  +57.41%  "
$ws.Range("E31").Value = "  -5.37%  "
$ws.Range("D32").Value = "4.44"
$ws.Range("E32").Value = "  -3.13%  "
$ws.Range("D33").Value = "0.0610"
$ws.Range("E33").Value = "  -1.77%  "
$ws.Range("D34").Value = "0.0926"
$ws.Range("E34").Value = "  +4.39%  "
$ws.Range("D35").Value = "0.973"
$ws.Range("E35").Value = "  +9.57%  "
$ws.Range("D36").Value = "2.39"
$ws.Range("E36").Value = "  +17.66%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("E38").Value = "  -1.96%  "
$ws.Range("D39").Value = "4.07"
$ws.Range("E39").Value = "  -6.29%  "
$ws.Range("E40").Value = "  -11.79%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "1.17"
$ws.Range("E41").Value = "  +5.19%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "0.0222"
$ws.Range("E42").Value = "  -2.48%  "
$ws.Range("D43").Value = "97.29"
$ws.Range("E43").Value = "  -8.64%  "
$ws.Range("D44").Value = "2.76"
$ws.Range("E44").Value = "  -6.27%  "
$ws.Range("D45").Value = "15.96"
$ws.Range("E45").Value = "  -9.31%  "
$ws.Range("D46").Value = "1.325.80"
$ws.Range("E46").Value = "  -1.56%  "
$ws.Range("D47").Value = "0.0841"
$ws.Range("E47").Value = "  +3.07%  "
$ws.Range("D48").Value = "6.96"
$ws.Range("E48").Value = "  +8.76%  "
$ws.Range("D49").Value = "2.294.29"
$ws.Range("E49").Value = "  +9.80%  "
$ws.Range("D50").Value = "2.85"
$ws.Range("E50").Value = "  +1.46%  "
$ws.Range("E51").Value = "  -6.49%  "
